$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 37625
$ws.Range("E2").Value = 736259302345
$ws.Range("F2").Value = 3517346364
$ws.Range("G2").Value = -0.22492

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 2079.34
$ws.Range("E3").Value = 250171234092
$ws.Range("F3").Value = 6903048294
$ws.Range("G3").Value = 0.12591

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 88749057889
$ws.Range("F4").Value = 16309226426
$ws.Range("G4").Value = -0.06523

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 233.52
$ws.Range("E5").Value = 35936092441
$ws.Range("F5").Value = 326624169
$ws.Range("G5").Value = -0.11017

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "XRP"
$ws.Range("D6").Value = 0.635868
$ws.Range("E6").Value = 34055784577
$ws.Range("F6").Value = 673962406
$ws.Range("G6").Value = 2.64056

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "USDC"
$ws.Range("D7").Value = 0.999526
$ws.Range("E7").Value = 24692866617
$ws.Range("F7").Value = 3454915825
$ws.Range("G7").Value = -0.07556

$ws.Range("B8").Value = "SOL"
$ws.Range("C8").Value = "Solana"
$ws.Range("D8").Value = 57.73
$ws.Range("E8").Value = 24447092537
$ws.Range("F8").Value = 750421402
$ws.Range("G8").Value = -0.54648

$ws.Range("B9").Value = "STETH"
$ws.Range("C9").Value = "Lido Staked Ether"
$ws.Range("D9").Value = 2080.02
$ws.Range("E9").Value = 19168159089
$ws.Range("F9").Value = 5566434
$ws.Range("G9").Value = 0.36862

$ws.Range("B10").Value = "ADA"
$ws.Range("C10").Value = "Cardano"
$ws.Range("D10").Value = 0.389958
$ws.Range("E10").Value = 13641644069
$ws.Range("F10").Value = 204086224
$ws.Range("G10").Value = -0.02658

$ws.Range("B11").Value = "DOGE"
$ws.Range("C11").Value = "Dogecoin"
$ws.Range("D11").Value = 0.077693
$ws.Range("E11").Value = 11028564613
$ws.Range("F11").Value = 320176930
$ws.Range("G11").Value = -0.57238

$ws.Range("B12").Value = "TRX"
$ws.Range("C12").Value = "TRON"
$ws.Range("D12").Value = 0.108234
$ws.Range("E12").Value = 9585065252
$ws.Range("F12").Value = 345169432
$ws.Range("G12").Value = 2.81656

$ws.Range("B13").Value = "LINK"
$ws.Range("C13").Value = "Chainlink"
$ws.Range("D13").Value = 15.08
$ws.Range("E13").Value = 8425946669
$ws.Range("F13").Value = 416451209
$ws.Range("G13").Value = 1.76866

$ws.Range("B14").Value = "TON"
$ws.Range("C14").Value = "Toncoin"
$ws.Range("D14").Value = 2.34
$ws.Range("E14").Value = 8158304662
$ws.Range("F14").Value = 28675464
$ws.Range("G14").Value = -1.68925

$ws.Range("B15").Value = "AVAX"
$ws.Range("C15").Value = "Avalanche"
$ws.Range("D15").Value = 21
$ws.Range("E15").Value = 7664721349
$ws.Range("F15").Value = 283746740
$ws.Range("G15").Value = 0.78611

$ws.Range("B16").Value = "MATIC"
$ws.Range("C16").Value = "Polygon"
$ws.Range("D16").Value = 0.772644
$ws.Range("E16").Value = 7163367529
$ws.Range("F16").Value = 279760699
$ws.Range("G16").Value = 0.05195

$ws.Range("B17").Value = "DOT"
$ws.Range("C17").Value = "Polkadot"
$ws.Range("D17").Value = 5.32
$ws.Range("E17").Value = 6929814098
$ws.Range("F17").Value = 125154645
$ws.Range("G17").Value = 0.53675

$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 37610
$ws.Range("E18").Value = 6074201967
$ws.Range("F18").Value = 69735164
$ws.Range("G18").Value = -0.16635

$ws.Range("B19").Value = "DAI"
$ws.Range("C19").Value = "Dai"
$ws.Range("D19").Value = 0.998215
$ws.Range("E19").Value = 5343350641
$ws.Range("F19").Value = 259647744
$ws.Range("G19").Value = -0.14416

$ws.Range("B20").Value = "LTC"
$ws.Range("C20").Value = "Litecoin"
$ws.Range("D20").Value = 70.76
$ws.Range("E20").Value = 5231368559
$ws.Range("F20").Value = 202438192
$ws.Range("G20").Value = -0.47644

$ws.Range("B21").Value = "SHIB"
$ws.Range("C21").Value = "Shiba Inu"
$ws.Range("D21").Value = 0.00000831
$ws.Range("E21").Value = 4902843612
$ws.Range("F21").Value = 114301831
$ws.Range("G21").Value = -0.02722

$ws.Range("B22").Value = "UNI"
$ws.Range("C22").Value = "Uniswap"
$ws.Range("D22").Value = 6.05
$ws.Range("E22").Value = 4574283751
$ws.Range("F22").Value = 189062114
$ws.Range("G22").Value = -1.89668

$ws.Range("B23").Value = "BCH"
$ws.Range("C23").Value = "Bitcoin Cash"
$ws.Range("D23").Value = 227.81
$ws.Range("E23").Value = 4460955545
$ws.Range("F23").Value = 74078147
$ws.Range("G23").Value = -0.0746

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "LEO Token"
$ws.Range("D24").Value = 3.98
$ws.Range("E24").Value = 3694153322
$ws.Range("F24").Value = 773897
$ws.Range("G24").Value = 0.50341

$ws.Range("B25").Value = "OKB"
$ws.Range("C25").Value = "OKB"
$ws.Range("D25").Value = 57.93
$ws.Range("E25").Value = 3480226996
$ws.Range("F25").Value = 5335786
$ws.Range("G25").Value = -0.69929

$ws.Range("B26").Value = "XLM"
$ws.Range("C26").Value = "Stellar"
$ws.Range("D26").Value = 0.122306
$ws.Range("E26").Value = 3424261647
$ws.Range("F26").Value = 64942610
$ws.Range("G26").Value = 0.9911

$ws.Range("B27").Value = "TUSD"
$ws.Range("C27").Value = "TrueUSD"
$ws.Range("D27").Value = 0.997935
$ws.Range("E27").Value = 3166077989
$ws.Range("F27").Value = 89159985
$ws.Range("G27").Value = -0.08446

$ws.Range("B28").Value = "XMR"
$ws.Range("C28").Value = "Monero"
$ws.Range("D28").Value = 170.93
$ws.Range("E28").Value = 3097072361
$ws.Range("F28").Value = 53908513
$ws.Range("G28").Value = 0.05869

$ws.Range("B29").Value = "KAS"
$ws.Range("C29").Value = "Kaspa"
$ws.Range("D29").Value = 0.132229
$ws.Range("E29").Value = 2877422719
$ws.Range("F29").Value = 38454528
$ws.Range("G29").Value = -4.55337

$ws.Range("B30").Value = "ATOM"
$ws.Range("C30").Value = "Cosmos Hub"
$ws.Range("D30").Value = 9.63
$ws.Range("E30").Value = 2816881054
$ws.Range("F30").Value = 180594648
$ws.Range("G30").Value = 7.14754

$ws.Range("B31").Value = "ETC"
$ws.Range("C31").Value = "Ethereum Classic"
$ws.Range("D31").Value = 19.3
$ws.Range("E31").Value = 2768898134
$ws.Range("F31").Value = 77080672
$ws.Range("G31").Value = 0.05147

$ws.Range("B32").Value = "CRO"
$ws.Range("C32").Value = "Cronos"
$ws.Range("D32").Value = 0.095499
$ws.Range("E32").Value = 2518642602
$ws.Range("F32").Value = 9855637
$ws.Range("G32").Value = -0.96205

$ws.Range("B33").Value = "LDO"
$ws.Range("C33").Value = "Lido DAO"
$ws.Range("D33").Value = 2.48
$ws.Range("E33").Value = 2210278660
$ws.Range("F33").Value = 54562650
$ws.Range("G33").Value = 0.74139

$ws.Range("B34").Value = "FIL"
$ws.Range("C34").Value = "Filecoin"
$ws.Range("D34").Value = 4.64
$ws.Range("E34").Value = 2197865478
$ws.Range("F34").Value = 99172801
$ws.Range("G34").Value = -0.59208

$ws.Range("B35").Value = "HBAR"
$ws.Range("C35").Value = "Hedera"
$ws.Range("D35").Value = 0.063441
$ws.Range("E35").Value = 2126752404
$ws.Range("F35").Value = 31632008
$ws.Range("G35").Value = 0.79944

$ws.Range("B36").Value = "ICP"
$ws.Range("C36").Value = "Internet Computer"
$ws.Range("D36").Value = 4.62
$ws.Range("E36").Value = 2077475492
$ws.Range("F36").Value = 45637801
$ws.Range("G36").Value = -0.78127

$ws.Range("B37").Value = "APT"
$ws.Range("C37").Value = "Aptos"
$ws.Range("D37").Value = 7.31
$ws.Range("E37").Value = 2025376564
$ws.Range("F37").Value = 62320791
$ws.Range("G37").Value = -0.59286

$ws.Range("B38").Value = "NEAR"
$ws.Range("C38").Value = "NEAR Protocol"
$ws.Range("D38").Value = 1.86
$ws.Range("E38").Value = 1867777797
$ws.Range("F38").Value = 126974564
$ws.Range("G38").Value = 1.20068

$ws.Range("B39").Value = "BUSD"
$ws.Range("C39").Value = "BUSD"
$ws.Range("D39").Value = 0.999615
$ws.Range("E39").Value = 1740659335
$ws.Range("F39").Value = 1163919679
$ws.Range("G39").Value = -0.06348

$ws.Range("B40").Value = "IMX"
$ws.Range("C40").Value = "Immutable"
$ws.Range("D40").Value = 1.36
$ws.Range("E40").Value = 1709477315
$ws.Range("F40").Value = 516895990
$ws.Range("G40").Value = -0.26468

$ws.Range("B41").Value = "VET"
$ws.Range("C41").Value = "VeChain"
$ws.Range("D41").Value = 0.02332331
$ws.Range("E41").Value = 1694355234
$ws.Range("F41").Value = 119369299
$ws.Range("G41").Value = 8.83242

$ws.Range("B42").Value = "RUNE"
$ws.Range("C42").Value = "THORChain"
$ws.Range("D42").Value = 5.34
$ws.Range("E42").Value = 1607650362
$ws.Range("F42").Value = 249071320
$ws.Range("G42").Value = 1.23576

$ws.Range("B43").Value = "MNT"
$ws.Range("C43").Value = "Mantle"
$ws.Range("D43").Value = 0.514692
$ws.Range("E43").Value = 1595857577
$ws.Range("F43").Value = 5899884
$ws.Range("G43").Value = 0.75026

$ws.Range("B44").Value = "OP"
$ws.Range("C44").Value = "Optimism"
$ws.Range("D44").Value = 1.81
$ws.Range("E44").Value = 1591982737
$ws.Range("F44").Value = 85305610
$ws.Range("G44").Value = 1.84431

$ws.Range("B45").Value = "TAO"
$ws.Range("C45").Value = "Bittensor"
$ws.Range("D45").Value = 267.08
$ws.Range("E45").Value = 1507931258
$ws.Range("F45").Value = 6805136
$ws.Range("G45").Value = -1.27712

$ws.Range("B46").Value = "QNT"
$ws.Range("C46").Value = "Quant"
$ws.Range("D46").Value = 103.35
$ws.Range("E46").Value = 1499940147
$ws.Range("F46").Value = 29965834
$ws.Range("G46").Value = 2.68999

$ws.Range("B47").Value = "AAVE"
$ws.Range("C47").Value = "Aave"
$ws.Range("D47").Value = 100.67
$ws.Range("E47").Value = 1476802660
$ws.Range("F47").Value = 131182539
$ws.Range("G47").Value = 2.71381

$ws.Range("B48").Value = "INJ"
$ws.Range("C48").Value = "Injective"
$ws.Range("D48").Value = 16.61
$ws.Range("E48").Value = 1398934277
$ws.Range("F48").Value = 92506417
$ws.Range("G48").Value = 1.91013

$ws.Range("B49").Value = "GRT"
$ws.Range("C49").Value = "The Graph"
$ws.Range("D49").Value = 0.147035
$ws.Range("E49").Value = 1371488248
$ws.Range("F49").Value = 56695825
$ws.Range("G49").Value = -3.99781

$ws.Range("B50").Value = "MKR"
$ws.Range("C50").Value = "Maker"
$ws.Range("D50").Value = 1464.8
$ws.Range("E50").Value = 1346583457
$ws.Range("F50").Value = 64790759
$ws.Range("G50").Value = 0.75963

$ws.Range("B51").Value = "ARB"
$ws.Range("C51").Value = "Arbitrum"
$ws.Range("D51").Value = 1.046
$ws.Range("E51").Value = 1335801236
$ws.Range("F51").Value = 125533010
$ws.Range("G51").Value = -1.14785
